$d = $word.ActiveDocument

$replacements = @(
    @{old="40×37="; new="46×58="},
    @{old="37×24="; new="74×34="},
    @{old="61×27="; new="29×75="},
    @{old="92×31="; new="77×37="},
    @{old="56×84="; new="11×78="},
    @{old="35×28="; new="76×41="},
    @{old="49×82="; new="66×54="},
    @{old="24×13="; new="48×37="},
    @{old="79×14="; new="92×97="},
    @{old="66×58="; new="25×77="},
    @{old="24×89="; new="65×29="},
    @{old="17×51="; new="82×16="},
    @{old="67×59="; new="64×49="},
    @{old="90×23="; new="96×35="},
    @{old="12×15="; new="11×46="},
    @{old="16×91="; new="88×97="},
    @{old="54×96="; new="96×18="},
    @{old="43×85="; new="59×82="},
    @{old="43×54="; new="97×42="},
    @{old="30×44="; new="53×32="},
    @{old="88×75="; new="69×70="},
    @{old="95×44="; new="26×82="},
    @{old="77×94="; new="94×46="},
    @{old="18×11="; new="50×91="},
    @{old="76×37="; new="95×30="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
